# The "About" sheet gets a new date stamp in C1 (next to the title in A1),
# recording when this data file was last touched/published.
# Serial date 44307 == April 21, 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$dateCell = $ws.Range("C1")
$dateCell.Value = (Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0)
$dateCell.NumberFormat = "m/d/yyyy"
